$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$win = $excel.ActiveWindow
$win.LargeScroll(16,0,0,0)  # scroll down by pages?
$ws.Range("I12").Select()
